$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 12.12601564876786
$arr[0,1] = 8.540681091977376
$arr[0,2] = 5.399918773833328
$arr[0,3] = 0
$arr[0,4] = 26.04993215201298
$arr[0,5] = 3.645707402750573
$arr[0,6] = 0
$arr[0,7] = 23.72435141585325
$arr[0,8] = 0
$arr[0,9] = 9.095403787623821
$arr[0,10] = 10.71897347368574
$arr[0,11] = 0
$arr[0,12] = 19.02750658731475
$arr[0,13] = 23.37924058404403
$ws.Range("B2:O2").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 11.8582662520356
$arr[0,1] = 8.534180673310509
$arr[0,2] = 5.344569355677074
$arr[0,3] = 0
$arr[0,4] = 26.07590100322824
$arr[0,5] = 3.647387742249168
$arr[0,6] = 0
$arr[0,7] = 23.80208550683385
$arr[0,8] = 0
$arr[0,9] = 8.911379957029517
$arr[0,10] = 10.69168628176257
$arr[0,11] = 0
$arr[0,12] = 19.08533027732752
$arr[0,13] = 23.44208854180393
$ws.Range("B3:O3").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 11.6929730747475
$arr[0,1] = 8.530342464852543
$arr[0,2] = 5.309725783333639
$arr[0,3] = 0
$arr[0,4] = 26.09817365588555
$arr[0,5] = 3.648474970010732
$arr[0,6] = 0
$arr[0,7] = 23.85406706283348
$arr[0,8] = 0
$arr[0,9] = 8.797866635037629
$arr[0,10] = 10.67705472575774
$arr[0,11] = 0
$arr[0,12] = 19.12249965832604
$arr[0,13] = 23.48528305634772
$ws.Range("B4:O4").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 11.62548878110133
$arr[0,1] = 8.52881741483127
$arr[0,2] = 5.295316896404458
$arr[0,3] = 0
$arr[0,4] = 26.1088396274013
$arr[0,5] = 3.648932020779797
$arr[0,6] = 0
$arr[0,7] = 23.87631837641757
$arr[0,8] = 0
$arr[0,9] = 8.751541372943255
$arr[0,10] = 10.67163061956656
$arr[0,11] = 0
$arr[0,12] = 19.13806648013272
$arr[0,13] = 23.50404159647803
$ws.Range("B5:O5").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 11.6142784051405
$arr[0,1] = 8.528566556299483
$arr[0,2] = 5.2929118180089
$arr[0,3] = 0
$arr[0,4] = 26.11070666062827
$arr[0,5] = 3.649008760327657
$arr[0,6] = 0
$arr[0,7] = 23.88007769195212
$arr[0,8] = 0
$arr[0,9] = 8.743846902077793
$arr[0,10] = 10.67076259031099
$arr[0,11] = 0
$arr[0,12] = 19.14067674169886
$arr[0,13] = 23.50722623114856
$ws.Range("B6:O6").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 11.69206333519729
$arr[0,1] = 8.530321738641454
$arr[0,2] = 5.309532301574213
$arr[0,3] = 0
$arr[0,4] = 26.09831106678681
$arr[0,5] = 3.648481077228334
$arr[0,6] = 0
$arr[0,7] = 23.85436282712811
$arr[0,8] = 0
$arr[0,9] = 8.797242065090421
$arr[0,10] = 10.67697938887554
$arr[0,11] = 0
$arr[0,12] = 19.12270789576146
$arr[0,13] = 23.48553136059603
$ws.Range("B7:O7").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 12.03393946139941
$arr[0,1] = 8.538408311807119
$arr[0,2] = 5.381016161668819
$arr[0,3] = 0
$arr[0,4] = 26.05757222708474
$arr[0,5] = 3.646275292010718
$arr[0,6] = 0
$arr[0,7] = 23.75027097713764
$arr[0,8] = 0
$arr[0,9] = 9.032098306585794
$arr[0,10] = 10.7091271255143
$arr[0,11] = 0
$arr[0,12] = 19.04709927837525
$arr[0,13] = 23.39995368033752
$ws.Range("B8:O8").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 12.69310248870298
$arr[0,1] = 8.555456690494587
$arr[0,2] = 5.514121389987786
$arr[0,3] = 0
$arr[0,4] = 26.02792826475925
$arr[0,5] = 3.642388119125829
$arr[0,6] = 0
$arr[0,7] = 23.57992936482525
$arr[0,8] = 0
$arr[0,9] = 9.48581238160453
$arr[0,10] = 10.78879106153572
$arr[0,11] = 0
$arr[0,12] = 18.91198864994529
$arr[0,13] = 23.26875150923599
$ws.Range("B9:O9").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.16525270014776
$arr[0,1] = 8.568676719844879
$arr[0,2] = 5.607264718310963
$arr[0,3] = 0
$arr[0,4] = 26.036789572048
$arr[0,5] = 3.639796732946881
$arr[0,6] = 0
$arr[0,7] = 23.47542805931665
$arr[0,8] = 0
$arr[0,9] = 9.811541181438969
$arr[0,10] = 10.85712536491602
$arr[0,11] = 0
$arr[0,12] = 18.82066354895876
$arr[0,13] = 23.19477145047492
$ws.Range("B10:O10").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.37639071782557
$arr[0,1] = 8.574834980802329
$arr[0,2] = 5.648560278119503
$arr[0,3] = 0
$arr[0,4] = 26.04746069227117
$arr[0,5] = 3.63867470428185
$arr[0,6] = 0
$arr[0,7] = 23.43238253746979
$arr[0,8] = 0
$arr[0,9] = 9.957395755215819
$arr[0,10] = 10.89026171468794
$arr[0,11] = 0
$arr[0,12] = 18.78082444691141
$arr[0,13] = 23.1660001351101
$ws.Range("B11:O11").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.45574113991214
$arr[0,1] = 8.577187120970374
$arr[0,2] = 5.664037561111805
$arr[0,3] = 0
$arr[0,4] = 26.05245386218657
$arr[0,5] = 3.638257946136615
$arr[0,6] = 0
$arr[0,7] = 23.41672927634448
$arr[0,8] = 0
$arr[0,9] = 10.01224129250412
$arr[0,10] = 10.90309701599411
$arr[0,11] = 0
$arr[0,12] = 18.76598234690071
$arr[0,13] = 23.15580844982901
$ws.Range("B12:O12").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.43867961874347
$arr[0,1] = 8.576679659753687
$arr[0,2] = 5.660711476892804
$arr[0,3] = 0
$arr[0,4] = 26.05133619207724
$arr[0,5] = 3.638347341502414
$arr[0,6] = 0
$arr[0,7] = 23.42007168709648
$arr[0,8] = 0
$arr[0,9] = 10.00044731415776
$arr[0,10] = 10.90032005274466
$arr[0,11] = 0
$arr[0,12] = 18.76916801944943
$arr[0,13] = 23.15797210968729
$ws.Range("B13:O13").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.38293144048634
$arr[0,1] = 8.575028092742935
$arr[0,2] = 5.649836851423659
$arr[0,3] = 0
$arr[0,4] = 26.04785241660699
$arr[0,5] = 3.638640254627215
$arr[0,6] = 0
$arr[0,7] = 23.43108175832922
$arr[0,8] = 0
$arr[0,9] = 9.961915964905542
$arr[0,10] = 10.89131197010913
$arr[0,11] = 0
$arr[0,12] = 18.77959849244787
$arr[0,13] = 23.16514755536895
$ws.Range("B14:O14").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.34870330525098
$arr[0,1] = 8.574019064504265
$arr[0,2] = 5.643154766608895
$arr[0,3] = 0
$arr[0,4] = 26.04584242368936
$arr[0,5] = 3.638820729960695
$arr[0,6] = 0
$arr[0,7] = 23.43791005428158
$arr[0,8] = 0
$arr[0,9] = 9.938262570907073
$arr[0,10] = 10.88583143935383
$arr[0,11] = 0
$arr[0,12] = 18.7860192154919
$arr[0,13] = 23.16963435993145
$ws.Range("B15:O15").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.15137412027222
$arr[0,1] = 8.568277115501445
$arr[0,2] = 5.604543812694817
$arr[0,3] = 0
$arr[0,4] = 26.03622560203955
$arr[0,5] = 3.639871199910109
$arr[0,6] = 0
$arr[0,7] = 23.47833168440638
$arr[0,8] = 0
$arr[0,9] = 9.801957959336281
$arr[0,10] = 10.85500045188771
$arr[0,11] = 0
$arr[0,12] = 18.82330134604064
$arr[0,13] = 23.19675008667114
$ws.Range("B16:O16").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.02932758482916
$arr[0,1] = 8.564791158216078
$arr[0,2] = 5.580577554368793
$arr[0,3] = 0
$arr[0,4] = 26.03202509385085
$arr[0,5] = 3.640530150798622
$arr[0,6] = 0
$arr[0,7] = 23.5042805191297
$arr[0,8] = 0
$arr[0,9] = 9.717706499944036
$arr[0,10] = 10.83660660472689
$arr[0,11] = 0
$arr[0,12] = 18.8466086719019
$arr[0,13] = 23.21463621502047
$ws.Range("B17:O17").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 12.95879087735382
$arr[0,1] = 8.562799757412959
$arr[0,2] = 5.566691868874154
$arr[0,3] = 0
$arr[0,4] = 26.03023427358283
$arr[0,5] = 3.640914511062959
$arr[0,6] = 0
$arr[0,7] = 23.51962840776685
$arr[0,8] = 0
$arr[0,9] = 9.66903175944249
$arr[0,10] = 10.82622045821197
$arr[0,11] = 0
$arr[0,12] = 18.86017499027146
$arr[0,13] = 23.22538334873938
$ws.Range("B18:O18").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 12.93485266614429
$arr[0,1] = 8.562127863161773
$arr[0,2] = 5.561973258129349
$arr[0,3] = 0
$arr[0,4] = 26.02973537687884
$arr[0,5] = 3.641045568754274
$arr[0,6] = 0
$arr[0,7] = 23.52489752202449
$arr[0,8] = 0
$arr[0,9] = 9.652515969647416
$arr[0,10] = 10.82273734742274
$arr[0,11] = 0
$arr[0,12] = 18.86479592260548
$arr[0,13] = 23.22910102250259
$ws.Range("B19:O19").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.04235529605976
$arr[0,1] = 8.565160838213865
$arr[0,2] = 5.583139291331871
$arr[0,3] = 0
$arr[0,4] = 26.03240755529083
$arr[0,5] = 3.640459451017587
$arr[0,6] = 0
$arr[0,7] = 23.50147445447436
$arr[0,8] = 0
$arr[0,9] = 9.726697910009696
$arr[0,10] = 10.83854468758364
$arr[0,11] = 0
$arr[0,12] = 18.84411095776933
$arr[0,13] = 23.2126846392118
$ws.Range("B20:O20").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.39932299620256
$arr[0,1] = 8.575512656407847
$arr[0,2] = 5.653035391978122
$arr[0,3] = 0
$arr[0,4] = 26.04884986682405
$arr[0,5] = 3.638553998610799
$arr[0,6] = 0
$arr[0,7] = 23.42783026159464
$arr[0,8] = 0
$arr[0,9] = 9.973244446825086
$arr[0,10] = 10.89395012823011
$arr[0,11] = 0
$arr[0,12] = 18.77652819346456
$arr[0,13] = 23.16302085394487
$ws.Range("B21:O21").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.62907487306238
$arr[0,1] = 8.582395268110298
$arr[0,2] = 5.697778778299559
$arr[0,3] = 0
$arr[0,4] = 26.06514452233622
$arr[0,5] = 3.637356043978902
$arr[0,6] = 0
$arr[0,7] = 23.38347205920883
$arr[0,8] = 0
$arr[0,9] = 10.1321024127076
$arr[0,10] = 10.93183148208388
$arr[0,11] = 0
$arr[0,12] = 18.73378127081249
$arr[0,13] = 23.13466305370548
$ws.Range("B22:O22").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.50680112308224
$arr[0,1] = 8.57871138178694
$arr[0,2] = 5.673986046956185
$arr[0,3] = 0
$arr[0,4] = 26.05594112135021
$arr[0,5] = 3.637991093428369
$arr[0,6] = 0
$arr[0,7] = 23.40680132391749
$arr[0,8] = 0
$arr[0,9] = 10.04754157445672
$arr[0,10] = 10.91146324266229
$arr[0,11] = 0
$arr[0,12] = 18.75646632432394
$arr[0,13] = 23.14942258378963
$ws.Range("B23:O23").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 13.03646661587326
$arr[0,1] = 8.56499366621328
$arr[0,2] = 5.581981463405107
$arr[0,3] = 0
$arr[0,4] = 26.03223270023759
$arr[0,5] = 3.640491397188952
$arr[0,6] = 0
$arr[0,7] = 23.50274173833725
$arr[0,8] = 0
$arr[0,9] = 9.722633629090973
$arr[0,10] = 10.83766789191184
$arr[0,11] = 0
$arr[0,12] = 18.84523965520501
$arr[0,13] = 23.21356550089725
$ws.Range("B24:O24").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 12.51654930215889
$arr[0,1] = 8.55072039229595
$arr[0,2] = 5.478905555048015
$arr[0,3] = 0
$arr[0,4] = 26.03056461148693
$arr[0,5] = 3.643393054714499
$arr[0,6] = 0
$arr[0,7] = 23.62238896703289
$arr[0,8] = 0
$arr[0,9] = 9.364168968158586
$arr[0,10] = 10.76549226868907
$arr[0,11] = 0
$arr[0,12] = 18.94713940700867
$arr[0,13] = 23.30031514354943
$ws.Range("B25:O25").Value = $arr
